$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New departure rows (21-25) for Friday, Jan 13, mirroring the existing
# table layout: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=TO, F=SHORT,
# G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank), L=DIFFERENCE, M=(blank)

$rows = @(
    @{ Row=21; A=20; C="3:00 PM";  D="LO3994"; E="Warsaw";      F="(WAW)"; G="LOT ";       H="E75S"; I="(SP-LIA)"; J="3:03 PM"; L="0 hours, 3 minutes" },
    @{ Row=22; A=21; C="3:55 PM";  D="FR2474"; E="London";      F="(STN)"; G="Ryanair ";   H="B38M"; I="(EI-HMW)"; J="4:00 PM"; L="0 hours, 5 minutes" },
    @{ Row=23; A=22; C="8:00 PM";  D="FR5107"; E="Dublin";      F="(DUB)"; G="Ryanair ";   H="B738"; I="(EI-EMF)"; J="8:03 PM"; L="0 hours, 3 minutes" },
    @{ Row=24; A=23; C="9:00 PM";  D="W95176"; E="London";      F="(LTN)"; G="Wizz Air ";  H="A321"; I="(G-WUKG)"; J="8:55 PM"; L="0 hours, -5 minutes" },
    @{ Row=25; A=24; C="9:50 PM";  D="FR9504"; E="Bristol";     F="(BRS)"; G="Ryanair ";   H="B738"; I="(EI-DWH)"; J="9:44 PM"; L="0 hours, -6 minutes" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "Friday, Jan 13"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 12).Value = $r.L
}
